$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "95.062.73"
Set-TextValue $ws "E2" "  -1.57%  "

Set-TextValue $ws "D3" "3.570.02"
Set-TextValue $ws "E3" "  -1.97%  "

Set-TextValue $ws "D4" "0.999"
Set-TextValue $ws "E4" "  -0.09%  "

Set-TextValue $ws "D5" "236.42"
Set-TextValue $ws "E5" "  -2.35%  "

Set-TextValue $ws "D6" "656.69"
Set-TextValue $ws "E6" "  +2.52%  "

Set-TextValue $ws "D7" "1.49"
Set-TextValue $ws "E7" "  -1.09%  "

Set-TextValue $ws "D8" "0.399"
Set-TextValue $ws "E8" "  -1.17%  "

Set-TextValue $ws "E9" "  +0.07%  "

Set-TextValue $ws "D10" "1.01"
Set-TextValue $ws "E10" "  -1.00%  "

Set-TextValue $ws "D11" "3.570.31"
Set-TextValue $ws "E11" "  -1.94%  "

Set-TextValue $ws "D12" "0.203"
Set-TextValue $ws "E12" "  +0.61%  "

Set-TextValue $ws "D13" "42.38"
Set-TextValue $ws "E13" "  -2.25%  "

Set-TextValue $ws "D14" "6.45"
Set-TextValue $ws "E14" "  +0.93%  "

Set-TextValue $ws "D15" "4.230.87"
Set-TextValue $ws "E15" "  -2.36%  "

Set-TextValue $ws "D16" "94.944.61"
Set-TextValue $ws "E16" "  -1.59%  "

Set-TextValue $ws "D17" "0.0000253"
Set-TextValue $ws "E17" "  -0.46%  "

Set-TextValue $ws "D18" "8.54"
Set-TextValue $ws "E18" "  +7.53%  "

Set-TextValue $ws "D19" "3.558.30"
Set-TextValue $ws "E19" "  -1.98%  "

Set-TextValue $ws "D20" "12.70"
Set-TextValue $ws "E20" "  -4.61%  "

Set-TextValue $ws "D21" "17.80"
Set-TextValue $ws "E21" "  -2.77%  "

Set-TextValue $ws "D22" "3.46"
Set-TextValue $ws "E22" "  -0.22%  "

Set-TextValue $ws "D23" "508.80"
Set-TextValue $ws "E23" "  -1.66%  "

Set-TextValue $ws "D24" "0.483"
Set-TextValue $ws "E24" "  -3.71%  "

Set-TextValue $ws "D25" "6.81"
Set-TextValue $ws "E25" "  +0.75%  "

Set-TextValue $ws "D26" "0.0000196"
Set-TextValue $ws "E26" "  -1.19%  "

Set-TextValue $ws "D27" "95.10"
Set-TextValue $ws "E27" "  -3.12%  "

Set-TextValue $ws "D28" "12.61"
Set-TextValue $ws "E28" "  +0.44%  "

Set-TextValue $ws "D29" "3.757.35"
Set-TextValue $ws "E29" "  -2.01%  "

Set-TextValue $ws "D30" "3.03"
Set-TextValue $ws "E30" "  -3.93%  "

Set-TextValue $ws "D31" "0.144"
Set-TextValue $ws "E31" "  -0.43%  "

Set-TextValue $ws "D32" "11.54"
Set-TextValue $ws "E32" "  -1.05%  "

Set-TextValue $ws "E33" "  +0.02%  "

Set-TextValue $ws "E34" "  +0.67%  "

Set-TextValue $ws "D35" "0.176"
Set-TextValue $ws "E35" "  -3.68%  "

Set-TextValue $ws "D36" "31.81"
Set-TextValue $ws "E36" "  +2.91%  "

Set-TextValue $ws "D37" "1.70"
Set-TextValue $ws "E37" "  +14.40%  "

Set-TextValue $ws "D38" "0.558"
Set-TextValue $ws "E38" "  -2.60%  "

Set-TextValue $ws "D39" "8.53"
Set-TextValue $ws "E39" "  +7.68%  "

Set-TextValue $ws "D40" "592.97"
Set-TextValue $ws "E40" "  +2.67%  "

Set-TextValue $ws "E41" "  +0.06%  "

Set-TextValue $ws "D42" "0.151"
Set-TextValue $ws "E42" "  -1.16%  "

Set-TextValue $ws "D43" "0.905"
Set-TextValue $ws "E43" "  -3.01%  "

Set-TextValue $ws "D44" "1.83"
Set-TextValue $ws "E44" "  +4.97%  "

Set-TextValue $ws "D45" "34.86"
Set-TextValue $ws "E45" "  +30.19%  "

Set-TextValue $ws "D46" "5.77"
Set-TextValue $ws "E46" "  +0.64%  "

Set-TextValue $ws "D47" "2.28"
Set-TextValue $ws "E47" "  +2.66%  "

Set-TextValue $ws "D48" "23.39"
Set-TextValue $ws "E48" "  -1.75%  "

Set-TextValue $ws "D49" "0.0414"
Set-TextValue $ws "E49" "  -4.30%  "

Set-TextValue $ws "D50" "3.54"
Set-TextValue $ws "E50" "  +0.23%  "

Set-TextValue $ws "D51" "8.18"
Set-TextValue $ws "E51" "  -0.48%  "
